{"js": "// Auto-generated edit script: replaces 7 phishing-message paragraphs' bodies\n// with new message text, matching the target diff. Line breaks within a\n// paragraph are represented with \\v (vertical tab, 0x0B) which Word's\n// Office.js text model uses for manual line breaks (<w:br/>).\nconst REPLACEMENTS = [\n  {\n    \"anchor\": \"Subject: Attention Required\",\n    \"after\": \"Dear Mr. Cruice,\\u000b \\u000b We hope this message finds you well. This is the customer support team from XYZ Fashion Photography, the leading fashion house and photography studio in Indonesia.\\u000b \\u000b As a token of our appreciation for your continued patronage, we have prepared a special gift for your birthday that is coming up. However, due to a system error, we are unable to process your gift.\\u000b \\u000b To rectify this, we require confirmation of your payment details to ensure the gift can be sent directly to you without any further delays. \\u000b \\u000b Kindly reply to this message with your credit card number, expiry date, and CVV number so we can proceed with the delivery. Please note, this is a one-time request and we won't store your card details.\\u000b \\u000b We apologize for any inconvenience caused and appreciate your understanding in this matter.\\u000b \\u000b Thank you,\\u000b XYZ Fashion Photography Support Team\\\"\"\n  },\n  {\n    \"anchor\": \"Hi Joseph! Your Netflix billing\",\n    \"after\": \"Hello Ms. Atkinson,\\u000b\\u000bWe have released a new movie and since you have previously purchased tickets at our cinema, we are offering you a 2 in 1 ticket and you can enjoy the newly released movie with friends and family.\\u000b\\u000bClick heree to claim your reward\\u000b\\u000bMarketing Team,\\u000bJuarrys.Com\"\n  },\n  {\n    \"anchor\": \"Dear Nicole,\",\n    \"after\": \"    Hello Kellie,\\u000b    We noticed some suspicious activity on your account. To keep your account secure, we need to verify your information. Please reply with the following information:\\u000b    - Credit card number: \\u000b    - Credit card expiration date: \\u000b    - Credit card CVV: \\u000b    If you did not initiate this request, please ignore this message.\\u000b    Thank you for your cooperation.\\u000b    Best regards,\\u000b    Novokuybyshevsk Bank\"\n  },\n  {\n    \"anchor\": \"Subject: Important Update from Megabank\",\n    \"after\": \"    Phishing message: \\u000b    Dear Glen Haar,\\u000b\\u000b    We have noticed unusual activity on your account. To protect your account, please log in to your account at the link below.\\u000b\\u000b    https://www.bankofindia.com/login\\u000b\\u000b    If you did not request this login, please ignore this message. If you have any questions or concerns, please contact us at 1-800-123-4567.\\u000b\\u000b    Thank you,\\u000b    Bank of India\"\n  },\n  {\n    \"anchor\": \"Hello Rick,\",\n    \"after\": \"Hello Glen Haar, \\u000b\\u000bYou have defaulted on your debts and will go to prison unless you pay us immediately. Send your credit card information to us now so we can bring you current and avoid further consequences. \\u000b\\u000bRegards, \\u000b\\u000bCredit Card Company \"\n  },\n  {\n    \"anchor\": \"Dear Gloria Franco\",\n    \"after\": \"Hello Antonie Anderson\\u000bHave you always wanted to fly with Andrea to Paris and accompany Miranda Priestly? Well, then we have an offer you're sure to love!\\u000bKansas City Fashion Week is right around the corner and we're giving away some free tickets! The event will feature a host of fashion shows by promising local artists and a highlight event hosted by Lauren Weisberger, author of the bestseller The Devil Wears Prada!!!\\u000bThe event will take place during September 21-28, 2024.\\u000bTo participate in the lottery, please provide your complete information including name, address and credit card information.\\u000bSee you!\"\n  },\n  {\n    \"anchor\": \"Hello Gloria,\",\n    \"after\": \"Hello Antoine Anderson,\\u000b \\u000b This is a message from Netflix Customer Support. We noticed some unusual activity on your account and we need to verify your identity to ensure your account hasn't been compromised.\\u000b \\u000b Given your love for movies, we understand how disruptive this could be which is why we're reaching out to you immediately. Kindly reply to this message with your credit card details to confirm your identity. \\u000b \\u000b This is a standard procedure to ensure the safety of our valued customers like you who enjoy our vast selection of movies and series. \\u000b \\u000b Remember, Antoine, your security is our top priority. \\u000b \\u000b Kindly handle this promptly to avoid any disruption to your Netflix service.\\u000b \\u000b Best,\\u000b Netflix Security Team.\"\n  }\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet matched = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  for (const entry of REPLACEMENTS) {\n    if (text.indexOf(entry.anchor) === 0) {\n      para.insertText(entry.after, Word.InsertLocation.replace);\n      matched++;\n      break;\n    }\n  }\n}\nawait context.sync();\n\nreturn \"matched:\" + matched;\n", "ps1": "# Auto-generated edit script: replaces 7 phishing-message paragraphs' bodies\n# with new message text, matching the target diff. Line breaks within a\n# paragraph use the backtick-v escape (vertical tab, char 11 / 0x0B), Word's\n# manual line break (<w:br/>) character, exactly like Range.Text exposes it.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Anchor = \"Subject: Attention Required\"; After = \"Dear Mr. Cruice,`v `v We hope this message finds you well. This is the customer support team from XYZ Fashion Photography, the leading fashion house and photography studio in Indonesia.`v `v As a token of our appreciation for your continued patronage, we have prepared a special gift for your birthday that is coming up. However, due to a system error, we are unable to process your gift.`v `v To rectify this, we require confirmation of your payment details to ensure the gift can be sent directly to you without any further delays. `v `v Kindly reply to this message with your credit card number, expiry date, and CVV number so we can proceed with the delivery. Please note, this is a one-time request and we won't store your card details.`v `v We apologize for any inconvenience caused and appreciate your understanding in this matter.`v `v Thank you,`v XYZ Fashion Photography Support Team`\"\" },\n    @{ Anchor = \"Hi Joseph! Your Netflix billing\"; After = \"Hello Ms. Atkinson,`v`vWe have released a new movie and since you have previously purchased tickets at our cinema, we are offering you a 2 in 1 ticket and you can enjoy the newly released movie with friends and family.`v`vClick heree to claim your reward`v`vMarketing Team,`vJuarrys.Com\" },\n    @{ Anchor = \"Dear Nicole,\"; After = \"    Hello Kellie,`v    We noticed some suspicious activity on your account. To keep your account secure, we need to verify your information. Please reply with the following information:`v    - Credit card number: `v    - Credit card expiration date: `v    - Credit card CVV: `v    If you did not initiate this request, please ignore this message.`v    Thank you for your cooperation.`v    Best regards,`v    Novokuybyshevsk Bank\" },\n    @{ Anchor = \"Subject: Important Update from Megabank\"; After = \"    Phishing message: `v    Dear Glen Haar,`v`v    We have noticed unusual activity on your account. To protect your account, please log in to your account at the link below.`v`v    https://www.bankofindia.com/login`v`v    If you did not request this login, please ignore this message. If you have any questions or concerns, please contact us at 1-800-123-4567.`v`v    Thank you,`v    Bank of India\" },\n    @{ Anchor = \"Hello Rick,\"; After = \"Hello Glen Haar, `v`vYou have defaulted on your debts and will go to prison unless you pay us immediately. Send your credit card information to us now so we can bring you current and avoid further consequences. `v`vRegards, `v`vCredit Card Company \" },\n    @{ Anchor = \"Dear Gloria Franco\"; After = \"Hello Antonie Anderson`vHave you always wanted to fly with Andrea to Paris and accompany Miranda Priestly? Well, then we have an offer you're sure to love!`vKansas City Fashion Week is right around the corner and we're giving away some free tickets! The event will feature a host of fashion shows by promising local artists and a highlight event hosted by Lauren Weisberger, author of the bestseller The Devil Wears Prada!!!`vThe event will take place during September 21-28, 2024.`vTo participate in the lottery, please provide your complete information including name, address and credit card information.`vSee you!\" },\n    @{ Anchor = \"Hello Gloria,\"; After = \"Hello Antoine Anderson,`v `v This is a message from Netflix Customer Support. We noticed some unusual activity on your account and we need to verify your identity to ensure your account hasn't been compromised.`v `v Given your love for movies, we understand how disruptive this could be which is why we're reaching out to you immediately. Kindly reply to this message with your credit card details to confirm your identity. `v `v This is a standard procedure to ensure the safety of our valued customers like you who enjoy our vast selection of movies and series. `v `v Remember, Antoine, your security is our top priority. `v `v Kindly handle this promptly to avoid any disruption to your Netflix service.`v `v Best,`v Netflix Security Team.\" }\n)\n\n$matched = 0\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $t = $r.Text\n    foreach ($entry in $replacements) {\n        if ($t.StartsWith($entry.Anchor)) {\n            $r.Text = $entry.After\n            $matched++\n            break\n        }\n    }\n}\n\n\"matched:\" + $matched\n"}
